$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, D (Fecha serial), M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$data = @(
    @(2, 44441, 160, 20000, 21000, 20500, 1025),
    @(3, 44333, 100, 19500, 20000, 19750, 988),
    @(4, 44435, 260, 20000, 22000, 21115, 1056),
    @(5, 44365, 100, 20000, 21000, 20500, 1025),
    @(6, 44776, 160, 23000, 24000, 23500, 1175),
    @(7, 44784, 160, 27000, 28000, 27500, 1375),
    @(8, 44326, 160, 19500, 20000, 19750, 988),
    @(9, 44428, 100, 20000, 21000, 20500, 1025),
    @(10, 44882, 120, 28000, 30000, 29000, 1450),
    @(11, 44410, 200, 20000, 21000, 20500, 1025),
    @(12, 44431, 160, 21000, 22000, 21500, 1075),
    @(13, 44434, 100, 20000, 21000, 20500, 1025),
    @(14, 44350, 160, 19000, 20000, 19500, 975),
    @(15, 44466, 100, 20000, 21000, 20500, 1025),
    @(16, 44417, 160, 20000, 21000, 20500, 1025),
    @(17, 44442, 140, 20000, 21000, 20500, 1025),
    @(18, 44336, 100, 19500, 20000, 19750, 988),
    @(19, 44364, 140, 20000, 21000, 20500, 1025),
    @(20, 44880, 100, 28000, 30000, 29000, 1450),
    @(21, 44343, 100, 19500, 20000, 19750, 988),
    @(22, 44427, 200, 20000, 21000, 20500, 1025),
    @(23, 44874, 240, 29000, 30000, 29500, 1475),
    @(24, 44473, 40, 19500, 20000, 19750, 988),
    @(25, 44448, 100, 20000, 21000, 20500, 1025),
    @(26, 44445, 160, 20000, 21000, 20500, 1025),
    @(27, 44809, 60, 27000, 28000, 27500, 1375),
    @(28, 44315, 100, 20000, 21000, 20500, 1025),
    @(29, 44782, 200, 23500, 24000, 23750, 1188),
    @(30, 44778, 100, 23000, 24000, 23500, 1175),
    @(31, 44462, 100, 19500, 20000, 19750, 988),
    @(32, 44420, 160, 20000, 21000, 20500, 1025),
    @(33, 44301, 100, 18000, 19000, 18500, 925),
    @(34, 44474, 200, 19000, 20000, 19500, 975),
    @(35, 44879, 100, 28000, 30000, 29000, 1450),
    @(36, 44781, 160, 23000, 24000, 23500, 1175),
    @(37, 44810, 100, 27000, 28000, 27500, 1375),
    @(38, 44418, 200, 20000, 21000, 20500, 1025),
    @(39, 44407, 160, 20000, 21000, 20500, 1025),
    @(40, 44467, 200, 20000, 21000, 20500, 1025),
    @(41, 44335, 200, 19000, 20000, 19500, 975)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value2 = $entry[1]   # D: Fecha
    $ws.Cells.Item($r, 13).Value2 = $entry[2]  # M: Volumen
    $ws.Cells.Item($r, 14).Value2 = $entry[3]  # N: Precio minimo
    $ws.Cells.Item($r, 15).Value2 = $entry[4]  # O: Precio maximo
    $ws.Cells.Item($r, 16).Value2 = $entry[5]  # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value2 = $entry[6]  # S: Precio $/Kg
}
